# Requisitos Não Funcionais - adicao do usuario comum
# Update non-functional requirement descriptions (RNF-02, RNF-04, RNF-05, RNF-06, RNF-08)
# and normalize cell alignment / wrap formatting to match the revised layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update requirement text (column B) ---
$ws.Range("B3").Value = 'Ser um sistema web'
$ws.Range("B5").Value = 'Utilizar as seguintes linguagens/tecnologias: Java + Spring Boot + Hibernate e Javascript + Node.js + Svelte'
$ws.Range("B6").Value = 'Utilizar MySQL como Sistema Gerenciador de Banco de Dados'
$ws.Range("B7").Value = 'Criptografar senhas'
$ws.Range("B9").Value = 'Utilizar "soft delete" em dados críticos: não apagar, mas inativar'

# --- Wrap the longer RNF-04 description and size the row to fit ---
$ws.Range("B5").WrapText = $true
$ws.Rows(5).RowHeight = 31.5

# --- Vertically center the identifier column (A2:A11), matching the rest of the table ---
$ws.Range("A2:A11").VerticalAlignment = -4108

# --- Restore the original selection / active cell ---
[void]$ws.Range("A12").Select()
